$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Análise")

$data = @(
    ,@("Ação do BB (BBAS3) voltou a estar atraente, alerta BTG","NEUTRAL",0.42)
    ,@("Estados e municípios contratam R`$ 19 bi do BB em crédito com garantia em 2024","NEUTRAL",0.55)
    ,@("BB mantém suspensão de consignado do INSS via correspondente","NEGATIVE",0.84)
    ,@("Dividendos de bancos: Itaú lidera pagamentos em 2024","POSITIVE",0.57)
    ,@("Banco do Brasil é eleito o mais sustentável do mundo pela sexta vez","POSITIVE",0.71)
    ,@("Ciclo de baixa do Ibovespa está próximo de terminar, diz PagBank","NEGATIVE",0.84)
    ,@("No Pain, No Gain: Veja 7 ações ‘blindadas’ para períodos ruins","NEUTRAL",0.56)
    ,@("Ibovespa se livra de crash do DeepSeek e dispara 2%","POSITIVE",0.78)
    ,@("Ibovespa: Veja os 10 destaques de hoje; Minerva cai mais de 6%","NEGATIVE",0.73)
    ,@("Ibovespa: Veja os 13 destaques de hoje; Minerva sobe mais de 5%","POSITIVE",0.67)
    ,@("Ibovespa: Veja os 11 destaques de hoje; Minerva cai mais de 5%","NEGATIVE",0.73)
    ,@("UBS BB eleva recomendação de Itaú (ITUB4) para “compra”","POSITIVE",0.48)
    ,@("Ibovespa: Veja os 13 destaques de hoje; Azul sobe mais de 6%","POSITIVE",0.66)
    ,@("Ação da Weg (WEGE3) é a que você precisa ter em fevereiro, diz Santander","NEUTRAL",0.74)
    ,@("BBB1 ao 25: considerando a inflação, qual premiação deu ao ganhador do reality o maior poder de compra?","NEUTRAL",0.61)
    ,@("Mesmo em cenário desafiador, é hora de investir, diz estrategista-chefe da XP","NEUTRAL",0.72)
    ,@("Realização de lucros: o que fez o Ibovespa cair 1,15% após 3 sessões seguidas de alta","NEGATIVE",0.79)
    ,@("As ações de dividendos mais indicadas para fevereiro: duas estatais lideram o ranking","NEUTRAL",0.69)
    ,@("Mesmo elevado, novo teto de juro do consignado do INSS pode ser insuficiente; entenda","NEGATIVE",0.78)
    ,@("BBA prefere qualidade a valor e rebaixa ação do Bradesco (BBDC4) para neutra","NEGATIVE",0.77)
    ,@("ETFs de Bitcoin nos EUA completam 1 ano: como foi e o que vem pela frente?","NEUTRAL",0.64)
    ,@("“Bottom fishing”: Morgan adota estratégia para Brasil e indica ações preferidas","NEUTRAL",0.64)
    ,@("BB assina termo com BC para devolver R`$ 20,6 mi a clientes por cobranças indevidas","NEGATIVE",0.82)
    ,@("Goldman vê Itaú posicionado para elevar rentabilidade e reforça ação como preferida","POSITIVE",0.79)
    ,@("Ibovespa sobe 4,86% em janeiro e tem melhor desempenho desde agosto; alta continuará?","POSITIVE",0.83)
    ,@("Renda fixa, ações ou internacional? Saiba tudo o que é preciso para investir em 2025","NEUTRAL",0.76)
    ,@("Ibovespa tem dia de forte ganho com sinais do Copom, apesar de alta de juros; entenda","POSITIVE",0.7)
    ,@("Ibovespa volta a superar 127 mil pontos, a caminho de 1ª alta mensal em cinco meses","POSITIVE",0.79)
    ,@("Ibovespa fecha em alta, acima dos 126 mil pontos, puxado por Vale e com Itaú em foco","POSITIVE",0.86)
    ,@("Por que o Ibovespa sobe cerca de 1% apesar do dia negativo dos mercados globais?","POSITIVE",0.44)
    ,@("O que esperar para bancões? Perspectivas de 2025 mais difícil deve ofuscar bom 4º tri","NEGATIVE",0.6)
    ,@("Ibovespa: por que dados de inflação dos EUA animam e fazem o índice subir mais de 1%?","NEGATIVE",0.54)
    ,@("Ao Vivo: Confira os vencedores do Prêmio Outliers InfoMoney","NEUTRAL",0.6)
    ,@("Temporada de balanços do 4T24 ganha destaque: veja ações e setores para ficar de olho","NEUTRAL",0.52)
    ,@("7 ações sobem mais de 20% e 10 caem mais de 5%: os destaques do Ibovespa em janeiro","NEGATIVE",0.34)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
